# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Sheet "Rules", cell B11 currently holds the text "R40" (a shared
# string). The commit replaces its content with the text "1", keeping the
# cell's existing style (s="23") untouched.
#
# A plain `Range.Value = "1"` assignment would be auto-coerced to the
# *number* 1 by Excel (since it is entirely digits), which would change
# the cell's type from string ("t=s") to numeric - not what the diff
# shows. To keep it a genuine text value we build it as a formula that
# evaluates to the text "1", then convert that formula to a static value
# via copy / paste-special-values, exactly as "Convert to values" works
# in the Excel UI. The cell keeps its original formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
